$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used with PasteSpecial to copy formatting only
# (reuses an existing cell style rather than synthesizing a brand-new one).
$xlPasteFormats = -4122

# Row 7: "就业单位类型" — formatted like B6 (style index 12)
$ws.Range("B7").Value = "就业单位类型"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)

# Row 8: "就业单位规模" — formatted like B27 (style index 11)
$ws.Range("B8").Value = "就业单位规模"
$ws.Range("B27").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)

# Row 9: "就业地区" — formatted like B6 (style index 12)
$ws.Range("B9").Value = "就业地区"
$ws.Range("B6").Copy()
$ws.Range("B9").PasteSpecial($xlPasteFormats)

# Row 10: "省内就业地区分布" — formatted like B6 (style index 12)
$ws.Range("B10").Value = "省内就业地区分布"
$ws.Range("B6").Copy()
$ws.Range("B10").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Move the active selection to B10 to match the saved view state
$ws.Range("B10").Select() | Out-Null
